{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2026-02-22 Sunday\", \"2026-02-23 Monday\"],\n  [\"674\u00f77=96, 2\", \"511\u00f73=170, 1\"],\n  [\"218\u00f74=54, 2\", \"396\u00f74=99, 0\"],\n  [\"777\u00f79=86, 3\", \"746\u00f78=93, 2\"],\n  [\"777\u00f74=194, 1\", \"504\u00f77=72, 0\"],\n  [\"748\u00f76=124, 4\", \"307\u00f75=61, 2\"],\n  [\"847\u00f74=211, 3\", \"613\u00f76=102, 1\"],\n  [\"428\u00f77=61, 1\", \"132\u00f75=26, 2\"],\n  [\"646\u00f72=323, 0\", \"496\u00f77=70, 6\"],\n  [\"188\u00f72=94, 0\", \"549\u00f73=183, 0\"],\n  [\"154\u00f76=25, 4\", \"486\u00f78=60, 6\"],\n  [\"242\u00f74=60, 2\", \"452\u00f79=50, 2\"],\n  [\"314\u00f76=52, 2\", \"744\u00f74=186, 0\"],\n  [\"306\u00f74=76, 2\", \"707\u00f72=353, 1\"],\n  [\"855\u00f76=142, 3\", \"422\u00f75=84, 2\"],\n  [\"329\u00f77=47, 0\", \"482\u00f78=60, 2\"],\n  [\"583\u00f78=72, 7\", \"596\u00f79=66, 2\"],\n  [\"849\u00f76=141, 3\", \"617\u00f79=68, 5\"],\n  [\"741\u00f79=82, 3\", \"249\u00f74=62, 1\"],\n  [\"662\u00f74=165, 2\", \"199\u00f74=49, 3\"],\n  [\"196\u00f76=32, 4\", \"486\u00f72=243, 0\"],\n  [\"880\u00f79=97, 7\", \"105\u00f77=15, 0\"],\n  [\"860\u00f72=430, 0\", \"619\u00f73=206, 1\"],\n  [\"666\u00f77=95, 1\", \"329\u00f74=82, 1\"],\n  [\"626\u00f77=89, 3\", \"296\u00f73=98, 2\"],\n  [\"364\u00f74=91, 0\", \"214\u00f79=23, 7\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  // Replace every occurrence found (each search string is unique in this\n  // document, but loop defensively in case of repeats).\n  for (let i = results.items.length - 1; i >= 0; i--) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $findText, $replaceText) {\n    $r = $doc.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #         ReplaceWith, Replace)\n    $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Text $d \"2026-02-22 Sunday\" \"2026-02-23 Monday\"\nReplace-Text $d \"674\u00f77=96, 2\" \"511\u00f73=170, 1\"\nReplace-Text $d \"218\u00f74=54, 2\" \"396\u00f74=99, 0\"\nReplace-Text $d \"777\u00f79=86, 3\" \"746\u00f78=93, 2\"\nReplace-Text $d \"777\u00f74=194, 1\" \"504\u00f77=72, 0\"\nReplace-Text $d \"748\u00f76=124, 4\" \"307\u00f75=61, 2\"\nReplace-Text $d \"847\u00f74=211, 3\" \"613\u00f76=102, 1\"\nReplace-Text $d \"428\u00f77=61, 1\" \"132\u00f75=26, 2\"\nReplace-Text $d \"646\u00f72=323, 0\" \"496\u00f77=70, 6\"\nReplace-Text $d \"188\u00f72=94, 0\" \"549\u00f73=183, 0\"\nReplace-Text $d \"154\u00f76=25, 4\" \"486\u00f78=60, 6\"\nReplace-Text $d \"242\u00f74=60, 2\" \"452\u00f79=50, 2\"\nReplace-Text $d \"314\u00f76=52, 2\" \"744\u00f74=186, 0\"\nReplace-Text $d \"306\u00f74=76, 2\" \"707\u00f72=353, 1\"\nReplace-Text $d \"855\u00f76=142, 3\" \"422\u00f75=84, 2\"\nReplace-Text $d \"329\u00f77=47, 0\" \"482\u00f78=60, 2\"\nReplace-Text $d \"583\u00f78=72, 7\" \"596\u00f79=66, 2\"\nReplace-Text $d \"849\u00f76=141, 3\" \"617\u00f79=68, 5\"\nReplace-Text $d \"741\u00f79=82, 3\" \"249\u00f74=62, 1\"\nReplace-Text $d \"662\u00f74=165, 2\" \"199\u00f74=49, 3\"\nReplace-Text $d \"196\u00f76=32, 4\" \"486\u00f72=243, 0\"\nReplace-Text $d \"880\u00f79=97, 7\" \"105\u00f77=15, 0\"\nReplace-Text $d \"860\u00f72=430, 0\" \"619\u00f73=206, 1\"\nReplace-Text $d \"666\u00f77=95, 1\" \"329\u00f74=82, 1\"\nReplace-Text $d \"626\u00f77=89, 3\" \"296\u00f73=98, 2\"\nReplace-Text $d \"364\u00f74=91, 0\" \"214\u00f79=23, 7\"\n"}
